# Natmi following Dr Hou advice
# Rebuild the LR-pair result table: the sending/target cluster set grew from
# {FAPs, sCs} to {ECs, FAPs, sCs}, so the 2-row table becomes a 6-row table
# (3 sending clusters x 2 target clusters) for the Bmp4 -> Bmpr1b pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..T
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity

$rows = @(
    @("ECs","Bmp4","Bmpr1b","FAPs",2,0.6666666666666666,5.269639666666666,15.808919,0.09922110188645328,0.09922110188645328,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,6.812254660674555,61.310291946071,0.08953999090312259,0.08953999090312259),
    @("ECs","Bmp4","Bmpr1b","sCs",2,0.6666666666666666,5.269639666666666,15.808919,0.09922110188645328,0.09922110188645328,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,0.7365445623962222,6.628901061566,0.009681110983330694,0.009681110983330694),
    @("FAPs","Bmp4","Bmpr1b","FAPs",3,1,39.012863,117.038589,0.7345662131494083,0.7345662131494083,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,50.43334546745567,453.900109207101,0.662893787638124,0.662893787638124),
    @("FAPs","Bmp4","Bmpr1b","sCs",3,1,39.012863,117.038589,0.7345662131494083,0.7345662131494083,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,5.452879878660667,49.075918907946,0.07167242551128429,0.07167242551128429),
    @("sCs","Bmp4","Bmpr1b","FAPs",3,1,8.827567333333334,26.482702,0.1662126849641383,0.1662126849641383,3,1,1.292736333333333,3.878209,0.9024289107934966,0.9024289107934966,11.41171702674644,102.705453240718,0.1499951322522499,0.1499951322522499),
    @("sCs","Bmp4","Bmpr1b","sCs",3,1,8.827567333333334,26.482702,0.1662126849641383,0.1662126849641383,2,0.6666666666666666,0.1397713333333333,0.419314,0.09757108920650338,0.09757108920650338,1.233840856269778,11.104567706428,0.01621755271188838,0.01621755271188838)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowValues = $rows[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}
